$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1).Range
$xml1 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:bidi w:val="0"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b/><w:bCs/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>ContosoLearn 价值主张</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$p1.InsertXML($xml1)

$p2 = $d.Paragraphs(2).Range
$xml2 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:bidi w:val="0"/></w:pPr><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b/><w:bCs/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>ContosoLearn</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t xml:space="preserve"> 不仅仅是一个 eLearning 平台，它是个性化学习和技能开发的</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b/><w:bCs/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>自适应 AI 助手</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>。</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>与其他平台不同，我们提供一个</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b/><w:bCs/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>简化但功能强大的分析系统</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>，提供可操作的见解，以增强学习体验，而不会让你不堪重负。</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>我们的平台旨在</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b/><w:bCs/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>适应学习风格</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>，根据进度和偏好提供</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b/><w:bCs/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>自定义建议</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>。</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$p2.InsertXML($xml2)

$p3 = $d.Paragraphs(3).Range
$xml3 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:bidi w:val="0"/></w:pPr><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>我们在对</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b/><w:bCs/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>原始内容</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>的承诺中脱颖而出。</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>我们的课程不仅是第三方信息的汇编，而且由领域的</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b/><w:bCs/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>专家精心策划和创建</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>，确保你获得独特而有价值的知识。</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>这一承诺包括提供关于您的进展、改进领域和下一步行动的清晰图景，而不会让您迷失在数据的海洋中。</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$p3.InsertXML($xml3)

$p4 = $d.Paragraphs(4).Range
$xml4 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:bidi w:val="0"/></w:pPr><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>面对 eLearning 市场的激烈竞争，ContosoLearn 通过不断创新和专注于</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b/><w:bCs/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>缩小传统教育的差距</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>而脱颖而出。</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>我们一直在寻找让学习更高效、更有效、更愉快的方法。</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>我们不仅关注你学到了什么，还关注你如何学习。</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$p4.InsertXML($xml4)

$p5 = $d.Paragraphs(5).Range
$xml5 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:bidi w:val="0"/></w:pPr><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>从本质上讲，ContosoLearn 不仅仅是一个 eLearning 平台。</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>这是一种个性化的学习体验，适合你，随你一起成长，并让你能够以最有效的方式学习。</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>有了 ContosoLearn，你不仅可以学习，还可以不断发展。</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="SimSun" w:eastAsia="SimSun" w:hAnsi="SimSun" w:cs="SimSun"/><w:b w:val="0"/><w:bCs w:val="0"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>我们致力于帮助你高效有效地掌握新技能，让学习成为一段愉快的旅程，而不是一项艰巨的任务。</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$p5.InsertXML($xml5)

Write-Host "OK"
